$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "320018799382"
    3  = "320018799393"
    4  = "320018799420"
    5  = "320018799441"
    6  = "320018799485"
    7  = "320018799500"
    8  = "320018799533"
    9  = "320018799625"
    10 = "320018799658"
    11 = "320018799670"
    12 = "320018799717"
    13 = "320018799739"
    14 = "320018799761"
    15 = "320018799783"
    16 = "320018799810"
    17 = "320018799831"
    18 = "320018799875"
    19 = "320018792701"
    20 = "320018792734"
    21 = "320018792756"
    22 = "320018792789"
    23 = "320018792790"
    24 = "320018792804"
    25 = "320018792815"
    26 = "320018792826"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 16).Value = $values[$row]
}
